# Adds two new columns, I ("I0") and J ("IF"), to the sheet,
# mirroring the header style used by the existing H ("IP") column,
# and fills in the corresponding data values for rows 2-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---
# Copy the format of the existing "IP" header (H1) onto the two new
# header cells so they pick up the same bold/bordered/centered style,
# then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (2-35) for columns I and J ---
$data = @{
    2  = @(5, 7)
    3  = @(9, 9)
    4  = @(7, 7)
    5  = @(9, 9)
    6  = @(6, 7)
    7  = @(7, 7)
    8  = @(7, 8)
    9  = @(8, 8)
    10 = @(3, 4)
    11 = @(7, 8)
    12 = @(9, 9)
    13 = @(6, 7)
    14 = @(6, 6)
    15 = @(7, 7)
    16 = @(8, 8)
    17 = @(8, 8)
    18 = @(6, 7)
    19 = @(6, 6)
    20 = @(5, 6)
    21 = @(5, 5)
    22 = @(7, 8)
    23 = @(6, 6)
    24 = @(5, 5)
    25 = @(6, 6)
    26 = @(8, 8)
    27 = @(7, 7)
    28 = @(7, 7)
    29 = @(5, 5)
    30 = @(7, 7)
    31 = @(5, 5)
    32 = @(7, 7)
    33 = @(4, 4)
    34 = @(1, 2)
    35 = @(2, 2)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
